$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.566240906715393
$ws.Range("B1").Value = 3.98604941368103
$ws.Range("C1").Value = 3.464866161346436
$ws.Range("D1").Value = 1.545438885688782
$ws.Range("E1").Value = 0.9526641964912415
